$wb = $excel.ActiveWorkbook

# Update Master sheet: replace employee placeholder names with real names
$master = $wb.Worksheets("Master")
$master.Range("A3").Value = "sairandhree sule"
$master.Range("A4").Value = "Ajay Wani"
$master.Range("A5").Value = "Akshay Patil"
$master.Range("A6").Value = "Pooja Joshi"
$master.Range("A7").Value = "Temp emp"

# Move the selection on the Master sheet to A8
$master.Activate()
$master.Range("A8").Select()

# Rename the per-employee tabs to have spaces / real names
$wb.Worksheets("Sheet1").Name = "f1 l1"
$wb.Worksheets("Sheet2").Name = "f2 l2"
$wb.Worksheets("emp1").Name = "sairandhree sule"
$wb.Worksheets("emp2").Name = "ajay wani"
$wb.Worksheets("emp3").Name = "akshay patil"
$wb.Worksheets("emp4").Name = "pooja joshi"
